$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is not number-like: assign directly as text.
$ws.Range("D2").Value = "34.560.30"
$ws.Range("E2").Value = "  +14.41%  "
$ws.Range("D3").Value = "1.796.87"
$ws.Range("E3").Value = "  +7.44%  "
$ws.Range("E4").Value = "  -0.70%  "
$ws.Range("E5").Value = "  +5.78%  "
$ws.Range("E6").Value = "  +5.49%  "
$ws.Range("E7").Value = "  -0.78%  "
$ws.Range("E8").Value = "  +6.13%  "
$ws.Range("E9").Value = "  +5.67%  "
$ws.Range("E10").Value = "  +6.79%  "
$ws.Range("E11").Value = "  +9.50%  "
$ws.Range("E12").Value = "  +1.81%  "
$ws.Range("D13").Value = "2.037.98"
$ws.Range("E13").Value = "  +6.53%  "
$ws.Range("D14").Value = "1.783.41"
$ws.Range("E14").Value = "  +6.64%  "
$ws.Range("E15").Value = "  +4.20%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "34.380.78"
$ws.Range("E16").Value = "  +13.69%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("E17").Value = "  +8.46%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("E18").Value = "  -4.89%  "
$ws.Range("E19").Value = "  +6.69%  "
$ws.Range("E20").Value = "  +7.30%  "
$ws.Range("D21").Value = "0.0₃0762"
$ws.Range("E21").Value = "  +6.16%  "
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("E23").Value = "  +4.62%  "
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("E25").Value = "  -2.08%  "
$ws.Range("E26").Value = "  +1.81%  "
$ws.Range("E27").Value = "  +6.32%  "
$ws.Range("E28").Value = "  +4.84%  "
$ws.Range("E29").Value = "  +6.06%  "
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("E31").Value = "  +10.60%  "
$ws.Range("E32").Value = "  +2.74%  "
$ws.Range("E33").Value = "  +5.87%  "
$ws.Range("E34").Value = "  +9.06%  "
$ws.Range("D35").Value = "1.575.51"
$ws.Range("E35").Value = "  +6.85%  "
$ws.Range("E36").Value = "  +6.28%  "
$ws.Range("E37").Value = "  +11.68%  "
$ws.Range("E38").Value = "  +3.04%  "
$ws.Range("E39").Value = "  +6.93%  "
$ws.Range("E40").Value = "  +5.10%  "
$ws.Range("E41").Value = "  +6.64%  "
$ws.Range("E42").Value = "  +2.42%  "
$ws.Range("E43").Value = "  +7.53%  "
$ws.Range("E44").Value = "  +5.23%  "
$ws.Range("E45").Value = "  +2.98%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("E46").Value = "  +2.73%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("E47").Value = "  +4.34%  "
$ws.Range("D48").Value = "1.931.45"
$ws.Range("E48").Value = "  +6.64%  "
$ws.Range("E49").Value = "  +5.44%  "
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("E51").Value = "  +21.96%  "

# Cells whose new value looks like a plain number (e.g. "0.991"): force
# text storage (matching the source inline-string cells) by pre-formatting
# as Text before assignment, then reverting the cell style so no stray
# number-format/style is left behind.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.991"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.11"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.551"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.990"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.58"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.31"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.283"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0676"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0923"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.643"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.34"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.18"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.05"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "264.75"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.996"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.50"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.41"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.17"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.57"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.87"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.117"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.16"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.992"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.83"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0513"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.20"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.58"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.84"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "88.92"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.630"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.84"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.922"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.12"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0519"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.05"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.56"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.75"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.994"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.40"

$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
